# filter price max fix / version
#
# Max price was set below the highest price for some items in NHL,
# preventing them from showing up in the standard infinite list.
# Bump the "max-price" column (G) on the apparel sheet from 150 to 400
# (over double the highest item price) and make "apparel" the active/
# selected sheet (it previously was "Special Circumstances").

$wb = $excel.ActiveWorkbook

$wsApparel = $wb.Worksheets.Item("apparel")

# Column G ("max-price") header is row 3; data rows are 4-16.
$wsApparel.Range("G4:G16").Value = 400

# "apparel" becomes the selected/active sheet, with G4:G16 selected and
# the view scrolled so column E is left-most - this also clears the
# previous selection (A1) tracked on "Special Circumstances".
$wsApparel.Activate()
$wsApparel.Range("G4:G16").Select()
$excel.ActiveWindow.ScrollColumn = 5

$wb.Save()
